$wb = $excel.ActiveWorkbook

# --- Step 1: Create the new "November 2019" sheet by copying "October 2019" ---
# Copying preserves styles, formulas, merged cells and page setup exactly.
$src = $wb.Worksheets.Item("October 2019")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)
$ws4 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4.Name = "November 2019"

# --- Step 2: Fix the header formula in E2 to reference "November 2019" ---
$ws4.Range("E2").Formula = '="    """&"November 2019"&""""&":"'

# --- Step 3: Replace the guild standings data (rows 4-53: rank, guild, contribution) ---
$arr = New-Object 'object[,]' 50,3
$arr[0,0] = 1; $arr[0,1] = "Eternal"; $arr[0,2] = 237814795
$arr[1,0] = 2; $arr[1,1] = "Smile"; $arr[1,2] = 222924716
$arr[2,0] = 3; $arr[2,1] = "Savages"; $arr[2,2] = 206424668
$arr[3,0] = 4; $arr[3,1] = "Elite"; $arr[3,2] = 202612882
$arr[4,0] = 5; $arr[4,1] = "Bounce"; $arr[4,2] = 195124955
$arr[5,0] = 6; $arr[5,1] = "Spring"; $arr[5,2] = 146174652
$arr[6,0] = 7; $arr[6,1] = "Sunset"; $arr[6,2] = 134839967
$arr[7,0] = 8; $arr[7,1] = "Epic"; $arr[7,2] = 119673264
$arr[8,0] = 9; $arr[8,1] = "Downtime"; $arr[8,2] = 113365977
$arr[9,0] = 10; $arr[9,1] = "Beaters"; $arr[9,2] = 112448493
$arr[10,0] = 11; $arr[10,1] = "RainSong"; $arr[10,2] = 99249305
$arr[11,0] = 12; $arr[11,1] = "Imperium"; $arr[11,2] = 96254842
$arr[12,0] = 13; $arr[12,1] = "lolicafe"; $arr[12,2] = 96210291
$arr[13,0] = 14; $arr[13,1] = "Gintama"; $arr[13,2] = 86205385
$arr[14,0] = 15; $arr[14,1] = "Undertale"; $arr[14,2] = 85678663
$arr[15,0] = 16; $arr[15,1] = "Remorse"; $arr[15,2] = 82758780
$arr[16,0] = 17; $arr[16,1] = "Maha"; $arr[16,2] = 82327678
$arr[17,0] = 18; $arr[17,1] = "Erda"; $arr[17,2] = 80545548
$arr[18,0] = 19; $arr[18,1] = "Atelier"; $arr[18,2] = 78049625
$arr[19,0] = 20; $arr[19,1] = "Sora"; $arr[19,2] = 75168538
$arr[20,0] = 21; $arr[20,1] = "Cleanse"; $arr[20,2] = 74658941
$arr[21,0] = 22; $arr[21,1] = "Broke"; $arr[21,2] = 73544475
$arr[22,0] = 23; $arr[22,1] = "Tama"; $arr[22,2] = 72353558
$arr[23,0] = 24; $arr[23,1] = "Oceania"; $arr[23,2] = 68755618
$arr[24,0] = 25; $arr[24,1] = "Lithe"; $arr[24,2] = 65982830
$arr[25,0] = 26; $arr[25,1] = "Revive"; $arr[25,2] = 65235043
$arr[26,0] = 27; $arr[26,1] = "Ravers"; $arr[26,2] = 62360160
$arr[27,0] = 28; $arr[27,1] = "Rising"; $arr[27,2] = 61904014
$arr[28,0] = 29; $arr[28,1] = "Sugar"; $arr[28,2] = 61854383
$arr[29,0] = 30; $arr[29,1] = "Artifacts"; $arr[29,2] = 61167779
$arr[30,0] = 31; $arr[30,1] = "Fabled"; $arr[30,2] = 55097872
$arr[31,0] = 32; $arr[31,1] = "Aloe"; $arr[31,2] = 53616057
$arr[32,0] = 33; $arr[32,1] = "Earnest"; $arr[32,2] = 50461914
$arr[33,0] = 34; $arr[33,1] = "Skyfall"; $arr[33,2] = 49672252
$arr[34,0] = 35; $arr[34,1] = "CyberThreat"; $arr[34,2] = 48991344
$arr[35,0] = 36; $arr[35,1] = "chigga"; $arr[35,2] = 48972339
$arr[36,0] = 37; $arr[36,1] = "Mystical"; $arr[36,2] = 47766075
$arr[37,0] = 38; $arr[37,1] = "Fandom"; $arr[37,2] = 46803556
$arr[38,0] = 39; $arr[38,1] = "Path"; $arr[38,2] = 45036671
$arr[39,0] = 40; $arr[39,1] = "Comity"; $arr[39,2] = 44448074
$arr[40,0] = 41; $arr[40,1] = "Howl"; $arr[40,2] = 43668039
$arr[41,0] = 42; $arr[41,1] = "Bubbles"; $arr[41,2] = 42657720
$arr[42,0] = 43; $arr[42,1] = "Coffee"; $arr[42,2] = 42494732
$arr[43,0] = 44; $arr[43,1] = "RainDrop"; $arr[43,2] = 42203249
$arr[44,0] = 45; $arr[44,1] = "Weibo"; $arr[44,2] = 41811578
$arr[45,0] = 46; $arr[45,1] = "Kingdom"; $arr[45,2] = 41143753
$arr[46,0] = 47; $arr[46,1] = "Reboot"; $arr[46,2] = 39936512
$arr[47,0] = 48; $arr[47,1] = "Exorcist"; $arr[47,2] = 39638954
$arr[48,0] = 49; $arr[48,1] = "Faction"; $arr[48,2] = 38974765
$arr[49,0] = 50; $arr[49,1] = "Prestigious"; $arr[49,2] = 36645409

$ws4.Range("B4:D53").Value2 = $arr

# --- Step 4: Update sheet view settings for all sheets ---

# August 2019: zoom change only
$ws1 = $wb.Worksheets.Item("August 2019")
$ws1.Activate()
$aw1 = $excel.ActiveWindow
$aw1.Zoom = 220
$aw1.ScrollRow = 27
$aw1.ScrollColumn = 1
$ws1.Range("C49").Select() | Out-Null

# September 2019: zoom change only
$ws2 = $wb.Worksheets.Item("September 2019")
$ws2.Activate()
$aw2 = $excel.ActiveWindow
$aw2.Zoom = 220
$aw2.ScrollRow = 19
$aw2.ScrollColumn = 2
$ws2.Range("C15").Select() | Out-Null

# October 2019: zoom change, no longer the selected tab, topLeftCell moves to A34
$ws3 = $wb.Worksheets.Item("October 2019")
$ws3.Activate()
$aw3 = $excel.ActiveWindow
$aw3.Zoom = 220
$aw3.ScrollRow = 34
$aw3.ScrollColumn = 1
$ws3.Range("E3").Select() | Out-Null

# November 2019: new active/selected sheet, top-left at A1, zoom 220
$ws4.Activate()
$aw4 = $excel.ActiveWindow
$aw4.Zoom = 220
$aw4.ScrollRow = 1
$aw4.ScrollColumn = 1
$ws4.Range("E2").Select() | Out-Null
